$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.808.14"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "2.289.51"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "97.23"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.79"
$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -2.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.42"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("E12").Value = "  -2.90%  "

$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.85"
$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").Value = "2.634.60"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("D17").Value = "2.288.35"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").Value = "43.792.17"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("E19").Value = "  +2.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.10"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("E22").Value = "  +9.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.82"
$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  -4.06%  "

$ws.Range("E25").Value = "  +5.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.35"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.46"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.69"
$ws.Range("E30").Value = "  -0.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.67"
$ws.Range("E31").Value = "  +1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.81"
$ws.Range("E32").Value = "  -2.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0896"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  +7.03%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -2.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.26"
$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.48"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.20"
$ws.Range("E45").Value = "  -3.22%  "

$ws.Range("E46").Value = "  -4.35%  "

$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.21"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.48"
$ws.Range("E49").Value = "  -2.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.441"
$ws.Range("E50").Value = "  +5.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.53"
$ws.Range("E51").Value = "  +10.95%  "
